$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Bugs related to Change Password." ->
#      "Bugs related to Change Password (DONE)."
#    with "(DONE)" and "." as two separate bold runs.
# ---------------------------------------------------------------------------
$searchRng = $d.Content
$found = $searchRng.Find.Execute("Bugs related to Change Password.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target paragraph 'Bugs related to Change Password.'"
}
$paraStart = $searchRng.Start
$periodPos = $searchRng.End - 1
$periodRange = $d.Range($periodPos, $periodPos + 1)
$periodRange.Text = " "

$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -eq $paraStart) {
        $para = $p
        break
    }
}
if ($null -eq $para) {
    throw "Could not re-locate target paragraph after edit"
}

# Insert a brand-new list paragraph right after it for the "Fetching data..."
# bullet (inherits the ListParagraph style + numbering from $para).
$para.Range.InsertParagraphAfter()
$newPara = $para.Next()
$newPara.Range.Text = "Fetching data based on users "

# -- Build bold "(DONE)" + "." (as two separate runs) in a scratch paragraph,
#    then transplant the FormattedText to the end of paragraph 1. Building it
#    in isolation (own paragraph, no neighboring runs) avoids a runtime quirk
#    where applying Font.BoldBi next to other runs leaks a stray <w:bCs/>
#    onto neighboring runs. --
$insPos1 = $para.Range.End - 1

$scratch1 = $d.Paragraphs.Add()
$scratch1.Range.Text = "(DONE)."
$p1Range = $scratch1.Range
$scratch1Index = $scratch1.Index

$s1a = $d.Range($p1Range.Start, $p1Range.End - 1)
$s1a.Find.Execute("(DONE)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s1a.Font.Bold = 1
$s1a.Font.BoldBi = 1

$s1b = $d.Range($p1Range.Start, $p1Range.End - 1)
$s1b.Find.Execute(".", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s1b.Font.Bold = 1
$s1b.Font.BoldBi = 1

$whole1 = $d.Range($p1Range.Start, $p1Range.End - 1)
$ft1 = $whole1.FormattedText

# Assign BEFORE deleting the scratch paragraph (FormattedText is a live
# reference, not a detached snapshot).
$ip1 = $d.Range($insPos1, $insPos1)
$ip1.FormattedText = $ft1

# Remove the scratch paragraph (positions shifted by the insert above, so
# re-fetch it by its paragraph index).
$scratch1b = $d.Paragraphs.Item($scratch1Index)
$r1b = $scratch1b.Range
$d.Range($r1b.Start, $r1b.End).Delete()

# ---------------------------------------------------------------------------
# 2) New paragraph: "Fetching data based on users (DONE)."
#    Here "(DONE)." is a single bold run (not split).
# ---------------------------------------------------------------------------
$insPos2 = $newPara.Range.End - 1

$scratch2 = $d.Paragraphs.Add()
$scratch2.Range.Text = "(DONE)."
$p2Range = $scratch2.Range
$scratch2Index = $scratch2.Index

$s2 = $d.Range($p2Range.Start, $p2Range.End - 1)
$s2.Find.Execute("(DONE).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s2.Font.Bold = 1
$s2.Font.BoldBi = 1

$whole2 = $d.Range($p2Range.Start, $p2Range.End - 1)
$ft2 = $whole2.FormattedText

$ip2 = $d.Range($insPos2, $insPos2)
$ip2.FormattedText = $ft2

$scratch2b = $d.Paragraphs.Item($scratch2Index)
$r2b = $scratch2b.Range
$d.Range($r2b.Start, $r2b.End).Delete()

Write-Output "Applied release-notes edit."
